# Documentazione + Diario + Media
# Update material descriptions in the "Analisi Costi" cost sheet and
# move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Update the shared-string text for the four materials that got more
# detailed descriptions (row/position stays the same).
$ws.Range("B8").Value  = "Switch ZyXEL 5 porte 10/100Mbps x2"
$ws.Range("B9").Value  = "Scheda Ethernet 100Mbps"
$ws.Range("B10").Value = "Server proxy Squid (Ubuntu 16.04)"
$ws.Range("B11").Value = "Laptop hp (W10)"

# Move the active selection from C14 to B12.
$ws.Range("B12").Select()

$wb.Save()
